# Applies the reference/documentation updates described in the commit
# "Added references and updated documentation" to the Indices.xlsx
# vignette workbook.
#
# Net effect on cell *content* (shared-string de-dup/re-ordering in the
# OOXML is irrelevant at the COM level -- only the following four cells
# on Sheet1 actually change their displayed text):
#   A9  -> appended " or Emergence index ($EI$)"
#   F9  -> inserted "@mockColdToleranceAdapted1972; " citation
#   A17 -> appended " (Allan, Vogel and Peterson; Erbach)" and reworded tail
#   F17 -> inserted "@erbachTillageContinuousCorn1982; " citation
# Plus the sheet's active selection moved from F29 to F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = 'Mean germination time or Mean length of incubation time ($\overline{T}$) or Germination resistance ($GR$) or Sprouting index ($SI$) or Emergence index ($EI$)'

$ws.Range("F9").Value = '[@edmond_effects_1958; @czabator_germination_1962; @smith_germinating_1964; @gordon_observations_1969; @gordon_germination_1971; @mockColdToleranceAdapted1972; @ellis_improved_1980 @labouriau_germinacao_1983; @ranal_how_2006]'

$ws.Range("A17").Value = 'Speed of germination or Germination rate Index or index of velocity of germination or Emergence rate index (Allan, Vogel and Peterson; Erbach) or Germination index (AOSA)'

$ws.Range("F17").Value = '[@throneberry_relation_1955; @maguire_speed_1962; @allan_seedling_1962; @kendrick_photocontrol_1969; @bouton_germination_1976; @erbachTillageContinuousCorn1982; @aosa_seed_1983; @khandakar_jute_1983; @bradbeer_seed_1988; @wardle_allelopathic_1991]'

# Update the saved view/selection (was topLeftCell=C28 / F29, now B8 / F9).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F9").Select() | Out-Null
